# "Generate Report for Handoff"
#
# The 0aa7b477-8e3e-4488-b6c6-334f45781283.md file got a newer handoff pass,
# so its handoff timestamps move forward on the Overview sheet and on each
# language sheet (zh-cn / de-de). Every other row is untouched.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: Latest Handoff Date for the 0aa7b477 row (row 5) -----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D5").Value = "2016-37-20 14:37:06"

# --- zh-cn sheet: Latest Handoff Datetime for the 0aa7b477 row (row 5) ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E5").Value = "2016-03-20 14:37:02"

# --- de-de sheet: Latest Handoff Datetime for the 0aa7b477 row (row 5) ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E5").Value = "2016-03-20 14:37:06"
